# Fruta / hortaliza, semanal
# Insert two new weekly price rows (366-367) for "Terminal Hortofrutícola Agro
# Chillán - Limón", pushing the existing data (previously rows 366-456) down
# to rows 368-458.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 366:456 down to 368:458, leaving two blank rows at 366:367
$ws.Rows("366:367").Insert()

# New row 366 - Limón, 1a amarillo
$row366 = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44508, 16, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "1a amarillo", 160, 6000, 6500, 6250, "`$/malla 16 kilos", "Región de O'Higgins", 391, 16)
for ($i = 0; $i -lt $row366.Length; $i++) {
    $ws.Cells.Item(366, $i + 1).Value = $row366[$i]
}

# New row 367 - Limón, 2a amarillo
$row367 = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44508, 16, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "2a amarillo", 120, 5000, 5500, 5250, "`$/malla 16 kilos", "Región de O'Higgins", 328, 16)
for ($i = 0; $i -lt $row367.Length; $i++) {
    $ws.Cells.Item(367, $i + 1).Value = $row367[$i]
}
